# A new article ("Trees, power lines downed by storm" / time bucket JSON)
# was inserted ahead of the existing "Storm leaves damage in its wake"
# entry, which pushed the historical-distance ranks of the two existing
# rows (2 and 3) past each other. Net effect: rows 2 and 3 swap their
# title / timestamp / historical-distance / uri values (the "time bucket"
# column, D, is identical for both rows so it is unaffected), while the
# hyperlink targets for column E stay attached to the same cells as
# before - only the displayed link text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes what row 3 used to be.
$ws.Range("A2").Value2 = "Trees, power lines downed by storm"
$ws.Range("B2").Value2 = "2009-06-17T00:00:00UTC"
$ws.Range("C2").Value2 = 167
$ws.Range("E2").Value2 = "http://glasgowdailytimes.com/local/x211927760/Trees-power-lines-downed-by-storm"

# Row 3 becomes what row 2 used to be.
$ws.Range("A3").Value2 = "Storm leaves damage in its wake"
$ws.Range("B3").Value2 = "2009-06-16T00:00:00UTC"
$ws.Range("C3").Value2 = 166
$ws.Range("E3").Value2 = "http://glasgowdailytimes.com/local/x211927706/Storm-leaves-damage-in-its-wake"
